# One last-minute edit to add the supervisor to the presentation frontpage.
# Also refreshes the cached "today" date/time fields (handout + notes
# masters) from 20/10/2024 -> 24/10/2024, reflecting the later save date.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Handout master footer date placeholder: 20/10/2024 -> 24/10/2024
# ---------------------------------------------------------------------
$handoutDate = $p.HandoutMaster.HeadersFooters.DateAndTime
$handoutDate.Text = "24/10/2024"

# ---------------------------------------------------------------------
# 2) Notes master footer date placeholder: 20/10/2024 -> 24/10/2024
# ---------------------------------------------------------------------
$notesDate = $p.NotesMaster.HeadersFooters.DateAndTime
$notesDate.Text = "24/10/2024"

# ---------------------------------------------------------------------
# 3) Front-page (slide 1) content placeholder: add the blank line,
#    the "Final: B.Sc: Informatik: Games Engineering" line and the
#    "Supervisor: ..." line below "München, 24. Oktober 2024".
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$contentShape = $slide1.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange

$tr.InsertAfter("`r`rFinal: B.Sc: Informatik: Games Engineering`rSupervisor: Prof. Dr. rer. nat. David Plecher")

# Split "Final: B.Sc: Informatik: Games Engineering" into the three runs
# that PowerPoint creates when "B.Sc" gets flagged by the spell-checker:
# "Final: " | "B.Sc" | ": Informatik: Games Engineering"
$finalPara = $contentShape.TextFrame.TextRange.Paragraphs(4, 1)
$bsc = $finalPara.Characters(8, 4)
$bsc.Text = $bsc.Text
